$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.26, -0.059, 0.066, 89, -65, 20),
    @(-0.125, 0.014, 0.261, 50, -4, -159),
    @(-0.206, -0.202, 0.443, -44, 8, 157),
    @(-0.196, -0.248, 0.418, 89, 29, 18),
    @(-0.125, -0.116, 0.074, -88, 56, 140),
    @(-0.013, -0.08699999999999999, 0.18, 122, 41, -83),
    @(0.277, 0.12, 0.514, -143, -9, -92),
    @(0.022, -0.058, 0.055, 125, -37, 146),
    @(0.07099999999999999, 0.156, 0.209, 100, -59, 134),
    @(-0.417, -0.094, 0.366, 63, -9, -99)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
